$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 11).Value = 6859
$ws.Cells.Item(3, 11).Value = 7101
$ws.Cells.Item(4, 11).Value = 1462
$ws.Cells.Item(5, 11).Value = 503
$ws.Cells.Item(6, 11).Value = 7775
$ws.Cells.Item(7, 11).Value = 23700

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 11).Value = 434
$ws.Cells.Item(3, 11).Value = 471
$ws.Cells.Item(6, 11).Value = 513
$ws.Cells.Item(7, 11).Value = 1552

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(3, 11).Value = 179
$ws.Cells.Item(6, 11).Value = 117
$ws.Cells.Item(7, 11).Value = 505

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(3, 11).Value = 362
$ws.Cells.Item(7, 11).Value = 1021

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 11).Value = 230
$ws.Cells.Item(3, 11).Value = 263
$ws.Cells.Item(7, 11).Value = 801

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(4, 11).Value = 20
$ws.Cells.Item(6, 11).Value = 202
$ws.Cells.Item(7, 11).Value = 550

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(2, 11).Value = 208
$ws.Cells.Item(6, 11).Value = 169
$ws.Cells.Item(7, 11).Value = 724
$ws.Cells.Item(8, 11).Value = 1552
$ws.Cells.Item(10, 11).Value = 135
$ws.Cells.Item(15, 11).Value = 248
$ws.Cells.Item(19, 11).Value = 696
$ws.Cells.Item(20, 11).Value = 574
$ws.Cells.Item(21, 11).Value = 77
$ws.Cells.Item(27, 11).Value = 223
$ws.Cells.Item(29, 11).Value = 1290
$ws.Cells.Item(31, 11).Value = 262
$ws.Cells.Item(33, 11).Value = 1021
$ws.Cells.Item(34, 11).Value = 136
$ws.Cells.Item(36, 11).Value = 297
$ws.Cells.Item(37, 11).Value = 801
$ws.Cells.Item(42, 11).Value = 874
$ws.Cells.Item(43, 11).Value = 192
$ws.Cells.Item(48, 11).Value = 308
$ws.Cells.Item(51, 11).Value = 297
$ws.Cells.Item(55, 11).Value = 253
$ws.Cells.Item(63, 11).Value = 62
$ws.Cells.Item(65, 11).Value = 550
$ws.Cells.Item(67, 11).Value = 917
$ws.Cells.Item(70, 11).Value = 42
$ws.Cells.Item(76, 11).Value = 319
$ws.Cells.Item(78, 11).Value = 278
$ws.Cells.Item(79, 11).Value = 586
$ws.Cells.Item(83, 11).Value = 505
$ws.Cells.Item(84, 11).Value = 191
$ws.Cells.Item(87, 11).Value = 47
$ws.Cells.Item(88, 11).Value = 258
$ws.Cells.Item(89, 11).Value = 355
$ws.Cells.Item(90, 11).Value = 224
$ws.Cells.Item(91, 11).Value = 281
$ws.Cells.Item(94, 11).Value = 317
$ws.Cells.Item(96, 11).Value = 254
$ws.Cells.Item(97, 11).Value = 186
$ws.Cells.Item(101, 11).Value = 23700

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(3, 11).Value = 67
$ws.Cells.Item(7, 11).Value = 262

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(3, 11).Value = 334
$ws.Cells.Item(7, 11).Value = 917

$ws = $wb.Worksheets.Item('South Deering')
$ws.Cells.Item(6, 11).Value = 36
$ws.Cells.Item(7, 11).Value = 191

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 11).Value = 364
$ws.Cells.Item(3, 11).Value = 460
$ws.Cells.Item(6, 11).Value = 375
$ws.Cells.Item(7, 11).Value = 1290

$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(2, 11).Value = 48
$ws.Cells.Item(6, 11).Value = 145
$ws.Cells.Item(7, 11).Value = 308

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(2, 11).Value = 205
$ws.Cells.Item(7, 11).Value = 696

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(4, 11).Value = 23
$ws.Cells.Item(6, 11).Value = 161
$ws.Cells.Item(7, 11).Value = 319

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Cells.Item(5, 11).Value = 6
$ws.Cells.Item(7, 11).Value = 169

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(2, 11).Value = 238
$ws.Cells.Item(6, 11).Value = 323
$ws.Cells.Item(7, 11).Value = 874

$ws = $wb.Worksheets.Item('Avondale')
$ws.Cells.Item(3, 11).Value = 23
$ws.Cells.Item(7, 11).Value = 135

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(2, 11).Value = 83
$ws.Cells.Item(4, 11).Value = 24
$ws.Cells.Item(7, 11).Value = 278

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(6, 11).Value = 86
$ws.Cells.Item(7, 11).Value = 253

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(2, 11).Value = 78
$ws.Cells.Item(7, 11).Value = 254

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Cells.Item(3, 11).Value = 133
$ws.Cells.Item(6, 11).Value = 57
$ws.Cells.Item(7, 11).Value = 281

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Cells.Item(6, 11).Value = 46
$ws.Cells.Item(7, 11).Value = 77

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(6, 11).Value = 145
$ws.Cells.Item(7, 11).Value = 586

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(2, 11).Value = 198
$ws.Cells.Item(3, 11).Value = 186
$ws.Cells.Item(4, 11).Value = 26
$ws.Cells.Item(7, 11).Value = 574

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(2, 11).Value = 115
$ws.Cells.Item(7, 11).Value = 297

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 11).Value = 241
$ws.Cells.Item(3, 11).Value = 231
$ws.Cells.Item(6, 11).Value = 197
$ws.Cells.Item(7, 11).Value = 724

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Cells.Item(6, 11).Value = 39
$ws.Cells.Item(7, 11).Value = 136

$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(6, 11).Value = 145
$ws.Cells.Item(7, 11).Value = 317

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(3, 11).Value = 63
$ws.Cells.Item(7, 11).Value = 248

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Cells.Item(6, 11).Value = 65
$ws.Cells.Item(7, 11).Value = 208

$ws = $wb.Worksheets.Item('West Town')
$ws.Cells.Item(3, 11).Value = 41
$ws.Cells.Item(6, 11).Value = 99
$ws.Cells.Item(7, 11).Value = 186

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Cells.Item(3, 11).Value = 12
$ws.Cells.Item(7, 11).Value = 42

$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(2, 11).Value = 66
$ws.Cells.Item(7, 11).Value = 258

$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(2, 11).Value = 100
$ws.Cells.Item(7, 11).Value = 355

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Cells.Item(3, 11).Value = 53
$ws.Cells.Item(7, 11).Value = 223

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(3, 11).Value = 64
$ws.Cells.Item(7, 11).Value = 224

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(2, 11).Value = 81
$ws.Cells.Item(3, 11).Value = 80
$ws.Cells.Item(4, 11).Value = 31
$ws.Cells.Item(7, 11).Value = 297

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Cells.Item(3, 11).Value = 54
$ws.Cells.Item(7, 11).Value = 192

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Cells.Item(3, 11).Value = 15
$ws.Cells.Item(7, 11).Value = 47
